# Auto-generated edit script applying the cryptos price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.39"
$ws.Range("E2").Value = "'-0.89%"
$ws.Range("D3").Value = "'26.98"
$ws.Range("E3").Value = "'-1.13%"
$ws.Range("D4").Value = "'4.637"
$ws.Range("E4").Value = "'-11.15%"
$ws.Range("D5").Value = "'0.05878"
$ws.Range("E5").Value = "'-0.81%"
$ws.Range("D6").Value = "'6.631"
$ws.Range("E6").Value = "'-1.11%"
$ws.Range("D7").Value = "'0.8591"
$ws.Range("E7").Value = "'-0.88%"
$ws.Range("D8").Value = "'0.9415"
$ws.Range("E8").Value = "'-6.37%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'-0.90%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.04299"
$ws.Range("E10").Value = "'20.73%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07091"
$ws.Range("E11").Value = "'-1.30%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03144"
$ws.Range("E12").Value = "'-0.20%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09164"
$ws.Range("E13").Value = "'-0.68%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001541"
$ws.Range("E14").Value = "'-0.81%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006063"
$ws.Range("E15").Value = "'0.29%"
$ws.Range("D16").Value = "'0.006221"
$ws.Range("E16").Value = "'5.23%"
$ws.Range("D17").Value = "'3.519"
$ws.Range("E18").Value = "'-1.83%"
$ws.Range("E19").Value = "'-0.57%"
$ws.Range("D20").Value = "'0.3054"
$ws.Range("E20").Value = "'-2.85%"
$ws.Range("E21").Value = "'-0.46%"
$ws.Range("D22").Value = "'3.823"
$ws.Range("E22").Value = "'8.55%"
$ws.Range("D23").Value = "'0.04234"
$ws.Range("E23").Value = "'0.73%"
$ws.Range("D24").Value = "'0.001224"
$ws.Range("E24").Value = "'0.58%"
$ws.Range("E25").Value = "'-5.23%"
$ws.Range("E26").Value = "'0.12%"
$ws.Range("E27").Value = "'30.61%"
$ws.Range("D41").Value = "'0.006271"
$ws.Range("E41").Value = "'-4.60%"
$ws.Range("D42").Value = "'0.1103"
$ws.Range("E42").Value = "'-0.08%"
$ws.Range("D43").Value = "'0.002430"
$ws.Range("E43").Value = "'10.58%"
$ws.Range("D44").Value = "'0.01144"
$ws.Range("E44").Value = "'5.57%"
$ws.Range("D45").Value = "'0.00005479"
$ws.Range("E45").Value = "'0.73%"
$ws.Range("E46").Value = "'0.15%"
$ws.Range("D47").Value = "'0.06751"
$ws.Range("E47").Value = "'-38.08%"
$ws.Range("D48").Value = "'0.2386"
$ws.Range("E48").Value = "'10,591.95%"
$ws.Range("E49").Value = "'0.15%"
$ws.Range("E50").Value = "'0.15%"
